$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster changes from ECs to FAPs, and downstream stats are recomputed
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 33.94639966666666
$ws.Range("N2").Value = 101.839199
$ws.Range("O2").Value = 0.9746097333921855
$ws.Range("P2").Value = 0.9746097333921855
$ws.Range("Q2").Value = 337.1098930580492
$ws.Range("R2").Value = 3033.989037522443
$ws.Range("S2").Value = 0.9427925588776398
$ws.Range("T2").Value = 0.9427925588776397

# Row 3: Target cluster changes from FAPs to MuSCs, and downstream stats are recomputed
$ws.Range("D3").Value = "MuSCs"
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8843623333333334
$ws.Range("N3").Value = 2.653087
$ws.Range("O3").Value = 0.02539026660781448
$ws.Range("P3").Value = 0.02539026660781448
$ws.Range("Q3").Value = 8.782294869028778
$ws.Range("R3").Value = 79.040653821259
$ws.Range("S3").Value = 0.0245613742666515
$ws.Range("T3").Value = 0.0245613742666515

# Row 4: Sending cluster changes from FAPs to MuSCs, Target cluster changes from MuSCs to FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.3351376666666667
$ws.Range("H4").Value = 1.005413
$ws.Range("I4").Value = 0.03264606685570879
$ws.Range("J4").Value = 0.03264606685570878
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 33.94639966666666
$ws.Range("N4").Value = 101.839199
$ws.Range("O4").Value = 0.9746097333921855
$ws.Range("P4").Value = 0.9746097333921855
$ws.Range("Q4").Value = 11.37671717602078
$ws.Range("R4").Value = 102.390454584187
$ws.Range("S4").Value = 0.03181717451454581
$ws.Range("T4").Value = 0.0318171745145458

# Row 5: Target cluster changes from ECs to MuSCs, and downstream stats are recomputed
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.3351376666666667
$ws.Range("I5").Value = 0.03264606685570879
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8843623333333334
$ws.Range("N5").Value = 2.653087
$ws.Range("O5").Value = 0.02539026660781448
$ws.Range("P5").Value = 0.02539026660781448
$ws.Range("Q5").Value = 0.2963831288812223
$ws.Range("R5").Value = 2.667448159931001
$ws.Range("S5").Value = 0.0008288923411629819
$ws.Range("T5").Value = 0.0008288923411629817

# Rows 6 and 7 (the ECs-related rows) are removed entirely, shrinking the table to a 2x2 cluster grid
$ws.Rows("6:7").Delete()
